# Add a "ManualCheck" worksheet (placed after "Tables") that reads the
# category prices from the "Tables" sheet and re-derives a per-category
# weighted monthly total, keyed by a "Year/Month" header row - used for a
# manual sanity check of the numbers before they're consumed elsewhere.

$wb = $excel.ActiveWorkbook
$tables = $wb.Worksheets.Item("Tables")

# New sheet goes right after "Tables" -> sheetId 3 / second tab.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tables)
$ws.Name = "ManualCheck"

# Header row: day-of-week columns reordered (Mon, Thu, Fri, Sun, Wed, Tue, Sat)
# plus a "Sum" label/total off in J/K.
$ws.Range("A1").Value = "Year/Month"
$ws.Range("B1").Value = "Monday"
$ws.Range("C1").Value = "Thursday"
$ws.Range("D1").Value = "Friday"
$ws.Range("E1").Value = "Sunday"
$ws.Range("F1").Value = "Wednesday"
$ws.Range("G1").Value = "Tuesday"
$ws.Range("H1").Value = "Saturday"
$ws.Range("J1").Value = "Sum"
$ws.Range("K1").Formula = "=SUM(I:I)"

# Row 2: the per-weekday price multipliers for this particular month.
$ws.Range("A2").Value = "2020, Apr"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 4

# Row 4: BM category prices (copied from "Tables", reordered to match the
# header above) plus the weighted total in column I.
$ws.Range("A4").Value = "BM"
$ws.Range("B4").Value = 1.5
$ws.Range("C4").Value = 1.5
$ws.Range("D4").Value = 1.5
$ws.Range("E4").Value = 1.5
$ws.Range("F4").Value = 1.5
$ws.Range("G4").Value = 1.5
$ws.Range("H4").Value = 1.5
$ws.Range("I4").Formula = "=B4*`$B`$2+C4*`$C`$2+D4*`$D`$2+E4*`$E`$2+F4*`$F`$2+G4*`$G`$2+H4*`$H`$2"

# Row 5: TOI category prices.
$ws.Range("A5").Value = "TOI"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 5
$ws.Range("I5").Formula = "=B5*`$B`$2+C5*`$C`$2+D5*`$D`$2+E5*`$E`$2+F5*`$F`$2+G5*`$G`$2+H5*`$H`$2"

# Row 6: HT category prices.
$ws.Range("A6").Value = "HT"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 4
$ws.Range("I6").Formula = "=B6*`$B`$2+C6*`$C`$2+D6*`$D`$2+E6*`$E`$2+F6*`$F`$2+G6*`$G`$2+H6*`$H`$2"

# Center-align everything that has content (mirrors the "Tables" sheet look).
$ws.Range("A1:H1").HorizontalAlignment = -4108
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("A2:H2").HorizontalAlignment = -4108
$ws.Range("A4:I4").HorizontalAlignment = -4108
$ws.Range("A5:I5").HorizontalAlignment = -4108
$ws.Range("A6:I6").HorizontalAlignment = -4108

# Leave the cursor parked on I4 in the new sheet, but keep "Tables" as the
# active/visible tab, with its original selection restored.
$ws.Range("I4").Select() | Out-Null
$tables.Activate() | Out-Null
$tables.Range("I1").Select() | Out-Null
